$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.615.47"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.838.50"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'314.66"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.4290"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'0.3680"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").Value = "'0.07287"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").Value = "'0.8722"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").Value = "'20.78"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "1.820.67"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").Value = "'5.439"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "'6.560"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "'0.06944"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'80.31"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'0.000008984"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'15.50"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "27.762.09"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "'5.169"
$ws.Range("E22").Value = "  +4.01%  "
$ws.Range("D23").Value = "'10.89"
$ws.Range("E23").Value = "  +4.10%  "
$ws.Range("D24").Value = "2.152.96"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("D25").Value = "'1.984"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "'153.60"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").Value = "'18.83"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "'5.243"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").Value = "'114.64"
$ws.Range("E29").Value = "  -5.40%  "
$ws.Range("D30").Value = "'1.848"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").Value = "'0.08885"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").Value = "'0.7716"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "'4.558"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.966"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.150"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("D36").Value = "'1.003"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.097"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05365"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").Value = "'0.01949"
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("D40").Value = "'2.823"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.1676"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.5110"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'6.643"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "'8.517"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").Value = "'10.51"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'106.38"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.06523"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.4726"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'1.002"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").Value = "'1.632"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "'1.783"
$ws.Range("E51").Value = "  +2.81%  "
